$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.053582725339496
$ws.Range("D2").Value = 1.050684185109144
$ws.Range("E2").Value = 1.059460595610241
$ws.Range("F2").Value = 1.068449830793772
$ws.Range("I2").Value = 1.042640248125781
$ws.Range("J2").Value = 1.0585993101037
$ws.Range("K2").Value = 1.053437778997844
$ws.Range("L2").Value = 1.062190038389462
$ws.Range("M2").Value = 1.071154981823424
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.055120979901498
$ws.Range("D3").Value = 1.051819178856494
$ws.Range("E3").Value = 1.060834737988353
$ws.Range("F3").Value = 1.069915583190494
$ws.Range("I3").Value = 1.043044498534006
$ws.Range("J3").Value = 1.059785649844038
$ws.Range("K3").Value = 1.054384437793381
$ws.Range("L3").Value = 1.063376998598946
$ws.Range("M3").Value = 1.072435100852805
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.056114749340793
$ws.Range("D4").Value = 1.052551975567683
$ws.Range("E4").Value = 1.061722642927601
$ws.Range("F4").Value = 1.070862850454584
$ws.Range("I4").Value = 1.043303724065596
$ws.Range("J4").Value = 1.060551317887194
$ws.Range("K4").Value = 1.054994790276159
$ws.Range("L4").Value = 1.064143252061468
$ws.Range("M4").Value = 1.0732617340235
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.056532158959303
$ws.Range("D5").Value = 1.052859660512452
$ws.Range("E5").Value = 1.062095623712309
$ws.Range("F5").Value = 1.071260807779602
$ws.Range("I5").Value = 1.04341214191551
$ws.Range("J5").Value = 1.060872738455014
$ws.Range("K5").Value = 1.055250860787709
$ws.Range("L5").Value = 1.06446496266182
$ws.Range("M5").Value = 1.073608852436558
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.056602222349224
$ws.Range("D6").Value = 1.052911299868099
$ws.Range("E6").Value = 1.062158231748392
$ws.Range("F6").Value = 1.07132761070696
$ws.Range("I6").Value = 1.043430312947022
$ws.Range("J6").Value = 1.06092667921604
$ws.Range("K6").Value = 1.055293825700704
$ws.Range("L6").Value = 1.064518954692997
$ws.Range("M6").Value = 1.073667111978363
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.056120328236416
$ws.Range("D7").Value = 1.052556088366841
$ws.Range("E7").Value = 1.061727627864075
$ws.Range("F7").Value = 1.070868169047317
$ws.Range("I7").Value = 1.043305174949912
$ws.Range("J7").Value = 1.060555614549202
$ws.Range("K7").Value = 1.054998213946794
$ws.Range("L7").Value = 1.064147552426453
$ws.Range("M7").Value = 1.073266373793515
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.054102918616167
$ws.Range("D8").Value = 1.051068100083999
$ws.Range("E8").Value = 1.059925257663561
$ws.Range("F8").Value = 1.068945435614732
$ws.Range("I8").Value = 1.042777354787207
$ws.Range("J8").Value = 1.059000651390739
$ws.Range("K8").Value = 1.053758164991217
$ws.Range("L8").Value = 1.062591551214005
$ws.Range("M8").Value = 1.071587957646325
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.050535484534558
$ws.Range("D9").Value = 1.048433437740554
$ws.Range("E9").Value = 1.056739331444527
$ws.Range("F9").Value = 1.065548043957477
$ws.Range("I9").Value = 1.04182915036453
$ws.Range("J9").Value = 1.05624522219175
$ws.Range("K9").Value = 1.05155598245278
$ws.Range("L9").Value = 1.059835706675254
$ws.Range("M9").Value = 1.068617152131452
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.048148273618456
$ws.Range("D10").Value = 1.046668177511426
$ws.Range("E10").Value = 1.054608300565137
$ws.Range("F10").Value = 1.063276435282976
$ws.Range("I10").Value = 1.041184679795718
$ws.Range("J10").Value = 1.054397539773978
$ws.Range("K10").Value = 1.05007608801135
$ws.Range("L10").Value = 1.057988705279499
$ws.Range("M10").Value = 1.066627327788519
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.047112353230527
$ws.Range("D11").Value = 1.045901636116947
$ws.Range("E11").Value = 1.053683769300382
$ws.Range("F11").Value = 1.062291121495491
$ws.Range("I11").Value = 1.040902656842014
$ws.Range("J11").Value = 1.05359484152151
$ws.Range("K11").Value = 1.049432416654926
$ws.Range("L11").Value = 1.05718653079822
$ws.Range("M11").Value = 1.065763419831231
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.046727219348552
$ws.Range("D12").Value = 1.045616575834351
$ws.Range("E12").Value = 1.053340081657253
$ws.Range("F12").Value = 1.061924869353787
$ws.Range("I12").Value = 1.040797452756536
$ws.Range("J12").Value = 1.053296279545141
$ws.Range("K12").Value = 1.049192891687034
$ws.Range("L12").Value = 1.056888197910704
$ws.Range("M12").Value = 1.065442172491162
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.046809847734835
$ws.Range("D13").Value = 1.045677737389597
$ws.Range("E13").Value = 1.053413816422767
$ws.Range("F13").Value = 1.062003443707594
$ws.Range("I13").Value = 1.040820039721496
$ws.Range("J13").Value = 1.0533603405109
$ws.Range("K13").Value = 1.049244290467046
$ws.Range("L13").Value = 1.056952208170827
$ws.Range("M13").Value = 1.065511097256437
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.047080525079264
$ws.Range("D14").Value = 1.045878079768541
$ws.Range("E14").Value = 1.053655365640672
$ws.Range("F14").Value = 1.062260852388176
$ws.Range("I14").Value = 1.040893969800478
$ws.Range("J14").Value = 1.05357017058287
$ws.Range("K14").Value = 1.049412626391095
$ws.Range("L14").Value = 1.057161878090717
$ws.Range("M14").Value = 1.065736872681568
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.047247252051797
$ws.Range("D15").Value = 1.046001473046265
$ws.Range("E15").Value = 1.05380415532975
$ws.Range("F15").Value = 1.062419415302402
$ws.Range("I15").Value = 1.040939461078956
$ws.Range("J15").Value = 1.05369940004098
$ws.Range("K15").Value = 1.049516285653828
$ws.Range("L15").Value = 1.057291013454849
$ws.Range("M15").Value = 1.065875933301978
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.048216976144536
$ws.Range("D16").Value = 1.046719004035604
$ws.Range("E16").Value = 1.054669620457449
$ws.Range("F16").Value = 1.063341790937362
$ws.Range("I16").Value = 1.041203334044455
$ws.Range("J16").Value = 1.054450755958531
$ws.Range("K16").Value = 1.050118745438875
$ws.Range("L16").Value = 1.058041891538119
$ws.Range("M16").Value = 1.066684613375063
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.048824651052474
$ws.Range("D17").Value = 1.047168506279302
$ws.Range("E17").Value = 1.055212021448083
$ws.Range("F17").Value = 1.063919913771201
$ws.Range("I17").Value = 1.041368059172599
$ws.Range("J17").Value = 1.054921349683436
$ws.Range("K17").Value = 1.050495880677664
$ws.Range("L17").Value = 1.058512246913747
$ws.Range("M17").Value = 1.067191255536798
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.049178881878058
$ws.Range("D18").Value = 1.04743048419277
$ws.Range("E18").Value = 1.055528223490304
$ws.Range("F18").Value = 1.064256960340069
$ws.Range("I18").Value = 1.041463854879971
$ws.Range("J18").Value = 1.055195584873936
$ws.Range("K18").Value = 1.050715581036973
$ws.Range("L18").Value = 1.058786365152878
$ws.Range("M18").Value = 1.067486549865137
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.04929962915696
$ws.Range("D19").Value = 1.047519776608264
$ws.Range("E19").Value = 1.055636011404707
$ws.Range("F19").Value = 1.064371857168913
$ws.Range("I19").Value = 1.041496470377313
$ws.Range("J19").Value = 1.055289049102698
$ws.Range("K19").Value = 1.0507904465061
$ws.Range("L19").Value = 1.058879793243854
$ws.Range("M19").Value = 1.067587200268405
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.048759475666046
$ws.Range("D20").Value = 1.04712030062099
$ws.Range("E20").Value = 1.0551538447293
$ws.Range("F20").Value = 1.063857903552679
$ws.Range("I20").Value = 1.041350415295923
$ws.Range("J20").Value = 1.054870885774762
$ws.Range("K20").Value = 1.050455446252405
$ws.Range("L20").Value = 1.058461806293427
$ws.Range("M20").Value = 1.067136920591993
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.047000827005694
$ws.Range("D21").Value = 1.045819093166184
$ws.Range("E21").Value = 1.053584243088155
$ws.Range("F21").Value = 1.062185059268652
$ws.Range("I21").Value = 1.040872211632162
$ws.Range("J21").Value = 1.0535083920727
$ws.Range("K21").Value = 1.049363067762377
$ws.Range("L21").Value = 1.05710014578674
$ws.Range("M21").Value = 1.065670397266158
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.045893082554855
$ws.Range("D22").Value = 1.044999044432368
$ws.Range("E22").Value = 1.052595773876913
$ws.Range("F22").Value = 1.061131751341958
$ws.Range("I22").Value = 1.040568951307832
$ws.Range("J22").Value = 1.052649396396173
$ws.Range("K22").Value = 1.048673716047515
$ws.Range("L22").Value = 1.056241873869333
$ws.Range("M22").Value = 1.064746286939809
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.046480512957785
$ws.Range("D23").Value = 1.045433952664997
$ws.Range("E23").Value = 1.053119934307371
$ws.Range("F23").Value = 1.061690277143906
$ws.Range("I23").Value = 1.040729962318451
$ws.Range("J23").Value = 1.053104990766367
$ws.Range("K23").Value = 1.049039396261854
$ws.Range("L23").Value = 1.056697065562581
$ws.Range("M23").Value = 1.065236372133717
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.048788926284317
$ws.Range("D24").Value = 1.047142083329463
$ws.Range("E24").Value = 1.055180132813806
$ws.Range("F24").Value = 1.063885923804207
$ws.Range("I24").Value = 1.041358388687136
$ws.Range("J24").Value = 1.054893689028454
$ws.Range("K24").Value = 1.050473717682464
$ws.Range("L24").Value = 1.058484598955912
$ws.Range("M24").Value = 1.067161472899284
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.051459286112721
$ws.Range("D25").Value = 1.049116091652548
$ws.Range("E25").Value = 1.057564187351187
$ws.Range("F25").Value = 1.066427496593095
$ws.Range("I25").Value = 1.04207644663955
$ws.Range("J25").Value = 1.0569594294197
$ws.Range("K25").Value = 1.052127352002485
$ws.Range("L25").Value = 1.060549852986283
$ws.Range("M25").Value = 1.069386784174056
